$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "URL"
$ws.Range("D1").Value = "Image"

$ws.Range("C2").Select()
